# fix(publipostage): Correct status name
#
# - shared string "bleu" -> "noir" (statut_label color name)
# - statut_name labels reworded from "résultat et / ou publication posté..."
#   to "résultat postés ou publiés..." (and matching "pas de résultat" label)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the longer / more specific strings first so that the shorter
# "résultat et / ou publication posté" match doesn't get applied as a
# partial match inside the "... dans les 12/36 mois" variants.
$null = $ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois")
$null = $ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois")
$null = $ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés")
$null = $ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés")

# Rename the "bleu" status colour to "noir" everywhere it's used.
$null = $ws.Cells.Replace("bleu", "noir")
